# mm gas: latest memory map
# Updates the BARREL sheet ("memory map") of the EPICS_datas workbook:
#  - removes the old duplicate H column (Input flow/Output flow/...) labels
#  - rebuilds / extends the J (offset) and K (name) columns down to row 36
#    with the up-to-date list of memory-map entries
#  - adds a new L column with "Only write" annotations for two rows
#  - applies border + centered formatting to the new J/K cells

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BARREL")

# ---------------------------------------------------------------
# 1. Remove the stale H column entries (now redundant / superseded)
# ---------------------------------------------------------------
$ws.Range("H5").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("H8").ClearContents()

# ---------------------------------------------------------------
# 2. Numeric "offset" rows 5-20: J holds the byte offset (number),
#    K holds the field name (text). Apply border + centered style.
# ---------------------------------------------------------------
$offsetRows = @(
    @{ Row = 5;  Offset = 0;  Name = "Threshold pressure warning (mbar)" },
    @{ Row = 6;  Offset = 4;  Name = "Threshold pressure fault (mbar)" },
    @{ Row = 7;  Offset = 8;  Name = "Threshold input flow warning (<) (L/h)" },
    @{ Row = 8;  Offset = 12; Name = "Threshols input flow fault (>) (L/h)" },
    @{ Row = 9;  Offset = 16; Name = "Threshold overpressure warning (>) (L/h)" },
    @{ Row = 10; Offset = 20; Name = "Threshold overpressure fault (% of input flow value)" },
    @{ Row = 11; Offset = 24; Name = "Threshold gap in/out warning (L/h)" },
    @{ Row = 12; Offset = 28; Name = "Threshold gap in/out fault (L/h)" },
    @{ Row = 13; Offset = 32; Name = "Time without fault on start system (mn)" },
    @{ Row = 14; Offset = 36; Name = "Time without fault on new setpoint (mn)" },
    @{ Row = 15; Offset = 40; Name = "Input flow setpoint (L/h)" },
    @{ Row = 16; Offset = 44; Name = "Time to input flow setpoint (second)" },
    @{ Row = 17; Offset = 48; Name = "Input flow" },
    @{ Row = 18; Offset = 52; Name = "Output flow" },
    @{ Row = 19; Offset = 56; Name = "Overpressure flow" },
    @{ Row = 20; Offset = 60; Name = "Input pressure" }
)

foreach ($item in $offsetRows) {
    $jCell = $ws.Cells.Item($item.Row, 10)
    $kCell = $ws.Cells.Item($item.Row, 11)

    $jCell.Borders.LineStyle = 1
    $jCell.HorizontalAlignment = -4108
    $jCell.Value = $item.Offset

    $kCell.Borders.LineStyle = 1
    $kCell.HorizontalAlignment = -4108
    $kCell.Value = $item.Name
}

# ---------------------------------------------------------------
# 3. Bit rows 21-28 (offset 65.0-65.7): J holds the "byte.bit" label
#    (stored as text), K holds the field name.
# ---------------------------------------------------------------
$bitRows65 = @(
    @{ Row = 21; Bit = "65.0"; Name = "Delta input/output flow fault" },
    @{ Row = 22; Bit = "65.1"; Name = "Input flow fault" },
    @{ Row = 23; Bit = "65.2"; Name = "Overpressure flow fault" },
    @{ Row = 24; Bit = "65.3"; Name = "Start gas command" },
    @{ Row = 25; Bit = "65.4"; Name = "Stop gas command" },
    @{ Row = 26; Bit = "65.5"; Name = "empty" },
    @{ Row = 27; Bit = "65.6"; Name = "empty" },
    @{ Row = 28; Bit = "65.7"; Name = "empty" }
)

foreach ($item in $bitRows65) {
    $jCell = $ws.Cells.Item($item.Row, 10)
    $kCell = $ws.Cells.Item($item.Row, 11)

    $jCell.Borders.LineStyle = 1
    $jCell.HorizontalAlignment = -4108
    $jCell.Value = "'" + $item.Bit

    $kCell.Borders.LineStyle = 1
    $kCell.HorizontalAlignment = -4108
    $kCell.Value = $item.Name
}
# rows 26-28's K cells keep the plain bordered style (no centering)
$ws.Range("K26").HorizontalAlignment = 1
$ws.Range("K27").HorizontalAlignment = 1
$ws.Range("K28").HorizontalAlignment = 1

# ---------------------------------------------------------------
# 4. Bit rows 29-31 (offset 64.0-64.2): J holds the "byte.bit" label
#    (stored as text), K holds the field name - plain border style
#    (no centering), matching the rest of column K.
# ---------------------------------------------------------------
$bitRows64a = @(
    @{ Row = 29; Bit = "64.0"; Name = "Barrel warning general" },
    @{ Row = 30; Bit = "64.1"; Name = "Barrel fault general" },
    @{ Row = 31; Bit = "64.2"; Name = "Barrel flow ok (input and output flow gas is correct)" }
)

foreach ($item in $bitRows64a) {
    $jCell = $ws.Cells.Item($item.Row, 10)
    $kCell = $ws.Cells.Item($item.Row, 11)

    $jCell.Borders.LineStyle = 1
    $jCell.HorizontalAlignment = -4108
    $jCell.Value = "'" + $item.Bit

    $kCell.Borders.LineStyle = 1
    $kCell.HorizontalAlignment = -4108
    $kCell.Value = $item.Name
}

# ---------------------------------------------------------------
# 5. Bit rows 32-36 (offset 64.3-64.7): J holds the "byte.bit" label
#    (stored as text), K holds the field name.
# ---------------------------------------------------------------
$bitRows64b = @(
    @{ Row = 32; Bit = "64.3"; Name = "Pressure warning" },
    @{ Row = 33; Bit = "64.4"; Name = "Input flow warning" },
    @{ Row = 34; Bit = "64.5"; Name = "Overpressure flow warning" },
    @{ Row = 35; Bit = "64.6"; Name = "Delta input/output flow warning" },
    @{ Row = 36; Bit = "64.7"; Name = "Pressure fault" }
)

foreach ($item in $bitRows64b) {
    $jCell = $ws.Cells.Item($item.Row, 10)
    $kCell = $ws.Cells.Item($item.Row, 11)

    $jCell.Borders.LineStyle = 1
    $jCell.HorizontalAlignment = -4108
    $jCell.Value = "'" + $item.Bit

    $kCell.Borders.LineStyle = 1
    $kCell.HorizontalAlignment = -4108
    $kCell.Value = $item.Name
}

# ---------------------------------------------------------------
# 6. New L column annotations ("Only write") for rows 32 and 33
# ---------------------------------------------------------------
$ws.Range("L32").Value = "Only write"
$ws.Range("L33").Value = "Only write"

# ---------------------------------------------------------------
# 7. Update the sheet view to match the latest selection/scroll position
# ---------------------------------------------------------------
$ws.Range("L14").Select()
